# Minor cleanup on the "Example: Constraint Checking" slide:
# the paragraph
#     if (!matchTypes(variable.type, expr))
# had its trailing ", " and "expr))" split across two separately
# formatted <a:r> runs (left over from an earlier edit). Re-merge them
# into the single run that carries the ", " run's formatting, matching
# how PowerPoint collapses adjacent runs with identical formatting
# once the text is touched again.

$p = $ppt.ActivePresentation

# Locate the shape that contains the target code line (search instead of
# hard-coding a slide/shape index, so this is resilient to reordering).
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame -eq -1) {
            if ($shape.TextFrame.TextRange.Text -like "*matchTypes(variable.type, expr))*") {
                $targetShape = $shape
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# Find the paragraph holding the "if (!matchTypes(...))" line, then find
# the two adjacent runs that split ", " and "expr))".
$runA = $null
$runB = $null
$paraCount = $tr.Paragraphs().Count
for ($pi = 1; $pi -le $paraCount; $pi++) {
    $para = $tr.Paragraphs($pi)
    if ($para.Text -like "*matchTypes(variable.type, expr))*") {
        $runCount = $para.Runs().Count
        $ri = 1
        while ($ri -le ($runCount + 10) -and $para.Runs($ri).Text -ne "") {
            $run = $para.Runs($ri)
            if ($run.Text -eq ", ") {
                $runA = $run
                $runB = $para.Runs($ri + 1)
                break
            }
            $ri = $ri + 1
        }
    }
}

# Re-set the combined character range spanning both runs. Setting text on
# a multi-run character range merges it back into a single run (using the
# first run's formatting), exactly mirroring the ", " + "expr))" -> ", expr))"
# collapse in the target edit.
$mergedText = $runA.Text + $runB.Text
$combined = $tr.Characters($runA.Start, $mergedText.Length)
$combined.Text = $mergedText
